# Applies the "coin sort pad" update:
#  1. Refreshes the datetimeFigureOut date placeholder text
#     (14.11.2024 -> 12.03.2025) on the slide master and every slide layout.
#  2. Rotates a set of EU country-code labels on slide 2's map.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholders (slide master + all slide layouts)
# ---------------------------------------------------------------------------
$oldDate = "14.11.2024"
$newDate = "12.03.2025"

$m = $p.SlideMaster
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $shp = $m.Shapes.Item($j)
    if ($shp.Name -like "Date*" -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $shp = $lay.Shapes.Item($j)
        if ($shp.Name -like "Date*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Country-code label rotation on slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$countryMap = @{
    "Textfeld 4"  = "HR"   # was NL
    "Textfeld 32" = "SK"   # was FI
    "Textfeld 37" = "NL"   # was HR
    "Textfeld 57" = "GR"   # was IE
    "Textfeld 60" = "PT"   # was GR
    "Textfeld 66" = "FI"   # was SI
    "Textfeld 69" = "SI"   # was PT
    "Textfeld 72" = "LU"   # was SK
    "Textfeld 75" = "IE"   # was LU
    "Textfeld 78" = "LT"   # was MT
    "Textfeld 90" = "EE"   # was LT
    "Textfeld 96" = "MT"   # was EE
}

for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($countryMap.ContainsKey($shp.Name)) {
        $shp.TextFrame.TextRange.Text = $countryMap[$shp.Name]
    }
}
